# Regenerate merged AHB files
#
# The header row of the "AHB-Diff" sheet carries the column labels for the
# two AHB format-versions being diffed ("_old" / "_new" suffixes). This
# regeneration renames those suffixes to the concrete format versions that
# were actually merged (FV2304 / FV2310), turns the sheet's data range into
# a proper Excel Table ("Table1") built from the (now renamed) header row,
# and freezes the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row -------------------------------------------
# Columns A:J were "<Feldname>_old"   -> "<Feldname>_FV2304"
# Column  K   stays "diff"
# Columns L:U were "<Feldname>_new"   -> "<Feldname>_FV2310"
$fieldNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $fieldNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fieldNames[$i] + "_FV2304"
}
$ws.Cells.Item(1, 11).Value = "diff"
for ($i = 0; $i -lt $fieldNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fieldNames[$i] + "_FV2310"
}

# --- 2. Turn A1:U81 into a native Excel Table ---------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U81"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row --------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
